$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at position 18 for L_DNK_1 (shifts existing row 18+ down by one)
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "L_DNK_1"
$ws.Range("B18").Value = "Q_BMWK"
$ws.Range("C18").Value = "Datenbank zum Deutschen Nachhaltigkeitskodex"
$ws.Range("D18").Value = "XXXDatenbank zum Deutschen Nachhaltigkeitskodex"
$ws.Range("E18").Value = "https://www.deutscher-nachhaltigkeitskodex.de/de/bericht/berichte-einsehen/"
$ws.Range("F18").Value = ""

# Insert new row at position 69 for L_IAB_1 (shifts existing row 69+ down by one)
$ws.Rows.Item(69).Insert()
$ws.Range("A69").Value = "L_IAB_1"
$ws.Range("B69").Value = "Q_IAB"
$ws.Range("C69").Value = "Daten zur Tarifbindung und betrieblichen Interessenvertretung"
$ws.Range("D69").Value = "XXXDaten zur Tarifbindung und betrieblichen Interessenvertretung"
$ws.Range("E69").Value = "https://iab.de/daten/daten-zur-tarifbindung-und-betrieblichen-interessenvertetung/"
$ws.Range("F69").Value = ""
